# Applies the "Upload new version with timestamp" edit:
#  - keeps existing item #42 row (row 48) but gives it a new name/price
#    (shared-string reshuffle in the source diff effectively renames it to
#    "شامبو الفيف 200 مل" at 85.00)
#  - inserts 4 new item rows (43-46) below it, pushing the totals / footer
#    rows down from 49/50 to 53/54
#  - updates the grand total and the generation timestamp in the footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert four fresh rows right after row 48 (old totals row 49 and
#    footer row 50 shift down to 53 / 54).
# ---------------------------------------------------------------------
$ws.Rows("49:52").Insert()

# Copy the formatting (number format / font / borders / alignment) of the
# template row (48) onto each of the 4 new rows, cell by cell, so every
# column keeps the same style as the rest of the item table.
for ($row = 49; $row -le 52; $row++) {
    for ($col = 1; $col -le 17; $col++) {
        $src = $ws.Cells.Item(48, $col)
        $dst = $ws.Cells.Item($row, $col)
        $src.Copy()
        $dst.PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------
# 2. Row heights, matching the authored workbook.
# ---------------------------------------------------------------------
$ws.Rows(49).RowHeight = 25.5
$ws.Rows(50).RowHeight = 24.75
$ws.Rows(51).RowHeight = 25.5
$ws.Rows(52).RowHeight = 25.5
$ws.Rows(53).RowHeight = 24.75

# ---------------------------------------------------------------------
# 3. Merge the label / quantity / price cells on the new rows just like
#    every other item row (A:B, C:G, H:K, L:M, N:O).
# ---------------------------------------------------------------------
foreach ($row in 49..52) {
    $ws.Range("A$row`:B$row").Merge()
    $ws.Range("C$row`:G$row").Merge()
    $ws.Range("H$row`:K$row").Merge()
    $ws.Range("L$row`:M$row").Merge()
    $ws.Range("N$row`:O$row").Merge()
}

# ---------------------------------------------------------------------
# 4. Cell values.
# ---------------------------------------------------------------------

# Row 48 - item 42 becomes "شامبو الفيف 200 مل" at 85.00 (qty 1:0)
$ws.Range("H48").Value = "1:0"
$ws.Range("N48").Value = "85.00"
$ws.Range("P48").Value = "85.0000"
$ws.Range("Q48").Value = "1:0"

# Row 49 - item 43 "شاور جل مود "
$ws.Range("A49").Value = 43
$ws.Range("C49").Value = "شاور جل مود "
$ws.Range("H49").Value = "6:0"
$ws.Range("L49").Value = "0"
$ws.Range("N49").Value = "85.00"
$ws.Range("P49").Value = "85.0000"
$ws.Range("Q49").Value = "1:0"

# Row 50 - item 44 "صوفي طويل جدا"
$ws.Range("A50").Value = 44
$ws.Range("C50").Value = "صوفي طويل جدا"
$ws.Range("H50").Value = "5:0"
$ws.Range("L50").Value = "0"
$ws.Range("N50").Value = "50.00"
$ws.Range("P50").Value = "50.0000"
$ws.Range("Q50").Value = "1:0"

# Row 51 - item 45 "فازلين بيور صغير " (the original item, now shifted down)
$ws.Range("A51").Value = 45
$ws.Range("C51").Value = "فازلين بيور صغير "
$ws.Range("H51").Value = "5:0"
$ws.Range("L51").Value = "0"
$ws.Range("N51").Value = "10.00"
$ws.Range("P51").Value = "10.0000"
$ws.Range("Q51").Value = "1:0"

# Row 52 - item 46 "فيانسيه كريم بخاخ"
$ws.Range("A52").Value = 46
$ws.Range("C52").Value = "فيانسيه كريم بخاخ"
$ws.Range("H52").Value = "4:0"
$ws.Range("L52").Value = "0"
$ws.Range("N52").Value = "55.00"
$ws.Range("P52").Value = "55.0000"
$ws.Range("Q52").Value = "1:0"

# ---------------------------------------------------------------------
# 5. Totals row (now row 53) and footer (now row 54).
# ---------------------------------------------------------------------
$ws.Range("P53").Value = 2397.025

$ws.Range("A54").Value = "Wednesday, 27 August, 2025 1:18 PM"
